# Set the "Правильный ответ" (correct answer) column (G) to the literal
# number 4 for every question row (2-24) - previously each row stored a
# shared-string/number copy of whichever option text/value was correct;
# now they all just point at option slot 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G24").Value = 4

# Update the on-screen selection/scroll position left behind by the edit:
# the sheet was scrolled right (so column C is the left-most visible
# column) and the active cell left on G25 (just below the data).
$ws.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G25").Select()
